$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.298623555000205
$ws.Range("C2").Value = 0.03441681897172089
$ws.Range("E2").Value = 0.467434553805262
$ws.Range("F2").Value = 2.252961650663323
$ws.Range("G2").Value = 0.002448853630257982
$ws.Range("J2").Value = 0.05505876934115239
$ws.Range("K2").Value = 0.2563543176156884
$ws.Range("M2").Value = 0.4047295761058578
$ws.Range("N2").Value = 1.677563516305566
$ws.Range("O2").Value = 2.903187944901831
$ws.Range("B3").Value = 0.268664072483574
$ws.Range("C3").Value = 0.03155665027905741
$ws.Range("E3").Value = 0.4583039965379641
$ws.Range("F3").Value = 2.240385623857648
$ws.Range("G3").Value = 0.002451054681535669
$ws.Range("J3").Value = 0.05541812041655358
$ws.Range("K3").Value = 0.226211499089942
$ws.Range("M3").Value = 0.3851145555615574
$ws.Range("N3").Value = 1.695186275403914
$ws.Range("O3").Value = 2.921717797552944
$ws.Range("B4").Value = 0.2503166588719523
$ws.Range("C4").Value = 0.02978792247699147
$ws.Range("E4").Value = 0.4529174375393197
$ws.Range("F4").Value = 2.233850835522588
$ws.Range("G4").Value = 0.00245247814346023
$ws.Range("J4").Value = 0.05565108186068102
$ws.Range("K4").Value = 0.2076986425946217
$ws.Range("M4").Value = 0.3732434814645487
$ws.Range("N4").Value = 1.706560772882264
$ws.Range("O4").Value = 2.93461226648931
$ws.Range("B5").Value = 0.2428523692134092
$ws.Range("C5").Value = 0.02906402440986255
$ws.Range("E5").Value = 0.4507776760589479
$ws.Range("F5").Value = 2.231486425928466
$ws.Range("G5").Value = 0.002453076375969118
$ws.Range("J5").Value = 0.05574911760701706
$ws.Range("K5").Value = 0.2001536162798914
$ws.Range("M5").Value = 0.3684495328059114
$ws.Range("N5").Value = 1.711335296381935
$ws.Range("O5").Value = 2.940248299149204
$ws.Range("B6").Value = 0.2416136922426517
$ws.Range("C6").Value = 0.02894363383146725
$ws.Range("E6").Value = 0.4504257128001967
$ws.Range("F6").Value = 2.231111853495946
$ws.Range("G6").Value = 0.002453176810156391
$ws.Range("J6").Value = 0.05576558385764585
$ws.Range("K6").Value = 0.1989007280775184
$ws.Range("M6").Value = 0.3676561411970312
$ws.Range("N6").Value = 1.712136517282824
$ws.Range("O6").Value = 2.941207194319617
$ws.Range("B7").Value = 0.2502159419867382
$ws.Range("C7").Value = 0.02977817234954472
$ws.Range("E7").Value = 0.4528883559406438
$ws.Range("F7").Value = 2.23381773923569
$ws.Range("G7").Value = 0.002452486137762424
$ws.Range("J7").Value = 0.05565239143848721
$ws.Range("K7").Value = 0.2075968908166317
$ws.Range("M7").Value = 0.373178651751978
$ws.Range("N7").Value = 1.706624599657731
$ws.Range("O7").Value = 2.934686731657891
$ws.Range("B8").Value = 0.2882838308075009
$ws.Range("C8").Value = 0.03343326204610264
$ws.Range("E8").Value = 0.4642407922999325
$ws.Range("F8").Value = 2.248379138391257
$ws.Range("G8").Value = 0.00244959763973098
$ws.Range("J8").Value = 0.05518011963284941
$ws.Range("K8").Value = 0.2459623402385773
$ws.Range("M8").Value = 0.3979306023690725
$ws.Range("N8").Value = 1.683524853164883
$ws.Range("O8").Value = 2.909262140833292
$ws.Range("B9").Value = 0.3633003814890401
$ws.Range("C9").Value = 0.04049993440804656
$ws.Range("E9").Value = 0.488244007607463
$ws.Range("F9").Value = 2.286351825669755
$ws.Range("G9").Value = 0.002444502213406256
$ws.Range("J9").Value = 0.05435152549193578
$ws.Range("K9").Value = 0.3211437410408848
$ws.Range("M9").Value = 0.4478327644070745
$ws.Range("N9").Value = 1.642621685712802
$ws.Range("O9").Value = 2.871445381947296
$ws.Range("B10").Value = 0.4186246218485508
$ws.Range("C10").Value = 0.04562915394060951
$ws.Range("E10").Value = 0.5069408832672906
$ws.Range("F10").Value = 2.319998355734953
$ws.Range("G10").Value = 0.002441102004141308
$ws.Range("J10").Value = 0.05380190815052632
$ws.Range("K10").Value = 0.3763348939390312
$ws.Range("M10").Value = 0.4853232010151061
$ws.Range("N10").Value = 1.615247596143996
$ws.Range("O10").Value = 2.851008854607102
$ws.Range("B11").Value = 0.4438360759431532
$ws.Range("C11").Value = 0.04794875257616127
$ws.Range("E11").Value = 0.5156773341620067
$ws.Range("F11").Value = 2.336555025449258
$ws.Range("G11").Value = 0.002439628999058967
$ws.Range("J11").Value = 0.0535646517222661
$ws.Range("K11").Value = 0.4014308166467231
$ws.Range("M11").Value = 0.5025576161731422
$ws.Range("N11").Value = 1.603375370641158
$ws.Range("O11").Value = 2.843308894242995
$ws.Range("B12").Value = 0.4533890140744177
$ws.Range("C12").Value = 0.04882512436472553
$ws.Range("E12").Value = 0.5190188035788594
$ws.Range("F12").Value = 2.343004474511858
$ws.Range("G12").Value = 0.002439081763396522
$ws.Range("J12").Value = 0.05347664018510301
$ws.Range("K12").Value = 0.4109321190989306
$ws.Range("M12").Value = 0.5091095540862653
$ws.Range("N12").Value = 1.598963133459497
$ws.Range("O12").Value = 2.840622846336316
$ws.Range("B13").Value = 0.4513313622732937
$ws.Range("C13").Value = 0.04863647202178356
$ws.Range("E13").Value = 0.5182976845037786
$ws.Range("F13").Value = 2.341607474524011
$ws.Range("G13").Value = 0.002439199151476954
$ws.Range("J13").Value = 0.05349551360852978
$ws.Range("K13").Value = 0.4088859383818431
$ws.Range("M13").Value = 0.5076973402377618
$ws.Range("N13").Value = 1.599909672283708
$ws.Range("O13").Value = 2.841191113718878
$ws.Range("B14").Value = 0.4446218860665851
$ws.Range("C14").Value = 0.04802089271304055
$ws.Range("E14").Value = 0.5159515746074845
$ws.Range("F14").Value = 2.337082022385232
$ws.Range("G14").Value = 0.002439583766392604
$ws.Range("J14").Value = 0.05355737425804907
$ws.Range("K14").Value = 0.4022125366073226
$ws.Range("M14").Value = 0.5030961354681978
$ws.Range("N14").Value = 1.603010699084702
$ws.Range("O14").Value = 2.843083305849916
$ws.Range("B15").Value = 0.4405128942843533
$ws.Range("C15").Value = 0.04764356972466999
$ws.Range("E15").Value = 0.5145188314330937
$ws.Range("F15").Value = 2.334333465723105
$ws.Range("G15").Value = 0.002439820727688836
$ws.Range("J15").Value = 0.05359550420858294
$ws.Range("K15").Value = 0.3981246150327706
$ws.Range("M15").Value = 0.5002810970966891
$ws.Range("N15").Value = 1.604921045692223
$ws.Range("O15").Value = 2.844272254273136
$ws.Range("B16").Value = 0.4169778449723935
$ws.Range("C16").Value = 0.04547728355825598
$ws.Range("E16").Value = 0.5063745806460389
$ws.Range("F16").Value = 2.318941501262657
$ws.Range("G16").Value = 0.002441199748817147
$ws.Range("J16").Value = 0.05381767001769511
$ws.Range("K16").Value = 0.3746945625717331
$ws.Range("M16").Value = 0.4842004912421203
$ws.Range("N16").Value = 1.616035157749907
$ws.Range("O16").Value = 2.851544208313186
$ws.Range("B17").Value = 0.4025508531624951
$ws.Range("C17").Value = 0.04414479849616271
$ws.Range("E17").Value = 0.5014374967557274
$ws.Range("F17").Value = 2.309819326365087
$ws.Range("G17").Value = 0.00244206459241548
$ws.Range("J17").Value = 0.05395722897891808
$ws.Range("K17").Value = 0.3603179136997028
$ws.Range("M17").Value = 0.4743814621153604
$ws.Range("N17").Value = 1.623001996616381
$ws.Range("O17").Value = 2.856414383682619
$ws.Range("B18").Value = 0.3942570140825978
$ws.Range("C18").Value = 0.04337710026767638
$ws.Range("E18").Value = 0.498619577880099
$ws.Range("F18").Value = 2.304690216945232
$ws.Range("G18").Value = 0.002442568973722969
$ws.Range("J18").Value = 0.05403870152900225
$ws.Range("K18").Value = 0.3520478534358062
$ws.Range("M18").Value = 0.4687507556649493
$ws.Range("N18").Value = 1.627063762985275
$ws.Range("O18").Value = 2.859365861181146
$ws.Range("B19").Value = 0.3914495950286891
$ws.Range("C19").Value = 0.04311695078756372
$ws.Range("E19").Value = 0.4976692187210077
$ws.Range("F19").Value = 2.302973809440374
$ws.Range("G19").Value = 0.002442740942954519
$ws.Range("J19").Value = 0.05406649327669211
$ws.Range("K19").Value = 0.3492475956717271
$ws.Range("M19").Value = 0.4668472112656019
$ws.Range("N19").Value = 1.628448387444825
$ws.Range("O19").Value = 2.860390987063283
$ws.Range("B20").Value = 0.4040862011599131
$ws.Range("C20").Value = 0.04428677739497289
$ws.Range("E20").Value = 0.5019608060962355
$ws.Range("F20").Value = 2.310778213796382
$ws.Range("G20").Value = 0.002441971809750189
$ws.Range("J20").Value = 0.0539422483203218
$ws.Range("K20").Value = 0.3618484382663496
$ws.Range("M20").Value = 0.4754249629569642
$ws.Range("N20").Value = 1.622254710836049
$ws.Range("O20").Value = 2.855880390519303
$ws.Range("B21").Value = 0.4465924650872921
$ws.Range("C21").Value = 0.04820175813125616
$ws.Range("E21").Value = 0.5166397844526784
$ws.Range("F21").Value = 2.338406378127416
$ws.Range("G21").Value = 0.002439470509614961
$ws.Range("J21").Value = 0.05353915458539316
$ws.Range("K21").Value = 0.4041727335611824
$ws.Range("M21").Value = 0.5044469267666756
$ws.Range("N21").Value = 1.602097585238026
$ws.Range("O21").Value = 2.842521286304446
$ws.Range("B22").Value = 0.4744069931331012
$ws.Range("C22").Value = 0.05074869146996264
$ws.Range("E22").Value = 0.526426622474375
$ws.Range("F22").Value = 2.357510954882585
$ws.Range("G22").Value = 0.002437897289722546
$ws.Range("J22").Value = 0.05328638806975849
$ws.Range("K22").Value = 0.4318223761812305
$ws.Range("M22").Value = 0.52356379868273
$ws.Range("N22").Value = 1.589410568236127
$ws.Range("O22").Value = 2.835129628275382
$ws.Range("B23").Value = 0.4595588714343251
$ws.Range("C23").Value = 0.04939043147732036
$ws.Range("E23").Value = 0.5211855423673484
$ws.Range("F23").Value = 2.347218609799597
$ws.Range("G23").Value = 0.002438731333071994
$ws.Range("J23").Value = 0.0534203184065607
$ws.Range("K23").Value = 0.4170664551641892
$ws.Range("M23").Value = 0.5133471724841385
$ws.Range("N23").Value = 1.596137303017574
$ws.Range("O23").Value = 2.838952099258762
$ws.Range("B24").Value = 0.4033920689558386
$ws.Range("C24").Value = 0.0442225938243439
$ws.Range("E24").Value = 0.5017241541558377
$ws.Range("F24").Value = 2.310344341473368
$ws.Range("G24").Value = 0.002442013734416874
$ws.Range("J24").Value = 0.05394901721871292
$ws.Range("K24").Value = 0.3611565028215864
$ws.Range("M24").Value = 0.4749531514052947
$ws.Range("N24").Value = 1.622592382937378
$ws.Range("O24").Value = 2.856121336769661
$ws.Range("B25").Value = 0.3429686544833714
$ws.Range("C25").Value = 0.03859913604028975
$ws.Range("E25").Value = 0.4815640737982321
$ws.Range("F25").Value = 2.275070279599689
$ws.Range("G25").Value = 0.002445820110854832
$ws.Range("J25").Value = 0.05456527271707845
$ws.Range("K25").Value = 0.3981246150327706
$ws.Range("M25").Value = 0.5002810970966891
$ws.Range("N25").Value = 1.604921045692223
$ws.Range("O25").Value = 2.844272254273136
